# "Update spacing issue on slide"
#
# The "CONVERGENT ZONE" textbox was too narrow (2309067 EMU / 181.8163 pt
# wide) and needed to be widened to 2572554 EMU / 202.5633 pt so its text
# no longer crowds/wraps. Two slides contain this textbox:
#   - Slide 5: widen the textbox AND bring it to the front of the
#     z-order (it ends up as the very last shape on the slide).
#   - Slide 6: widen the analogous textbox in place (no re-ordering).

$EMU_PER_POINT = 12700
$newWidthPt = 2572554 / $EMU_PER_POINT

$p = $ppt.ActivePresentation

function Find-ConvergentZoneShape($slide) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $text = $shp.TextFrame.TextRange.Text
            if ($text -like "*CONVERGENT*" -and $text -like "*ZONE*") {
                return $shp
            }
        }
    }
    return $null
}

# --- Slide 5: widen + bring to front ------------------------------------
$slide5 = $p.Slides.Item(5)
$shape5 = Find-ConvergentZoneShape $slide5
if ($shape5 -eq $null) {
    throw "Could not find the CONVERGENT/ZONE textbox on slide 5"
}

$shape5.Width = $newWidthPt
$shape5.ZOrder(0)   # msoBringToFront -> becomes the last shape on the slide

# --- Slide 6: widen only --------------------------------------------------
$slide6 = $p.Slides.Item(6)
$shape6 = Find-ConvergentZoneShape $slide6
if ($shape6 -eq $null) {
    throw "Could not find the CONVERGENT/ZONE textbox on slide 6"
}

$shape6.Width = $newWidthPt
